$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C13 value from 0 to 5 (this will trigger recalculation of the dependent formulas I13 and I19)
$ws.Range("C13").Value = 5

# Update selection to F14 to match the final state
$ws.Range("F14").Select()
